$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-104 down to 35-105
$ws.Rows("34:34").Insert()

# Populate the new row 34 with the new weekly record
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 45099
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = 100112012
$ws.Range("G34").Value = "Espinaca"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 250
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = 950
$ws.Range("N34").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 317
$ws.Range("Q34").Value = 3
$ws.Range("R34").Value = "Hortaliza"
